{"js": "// 1) \"...term\uc758 \ubcc0\ub3d9\" + <hidden _GoBack bookmark> + \"\uc131\uc5d0 max(v, 0)...\"\n//    becomes a single, unsplit run \"...term\uc758 \ubcc0\ub3d9\uc131\uc5d0  max(v, 0)...\".\n//    The search range spans the old hidden bookmark transparently, so\n//    replacing it with the same visible text collapses it back into one\n//    run (using the formatting of the matched range) and drops the\n//    bookmark that used to live inside it.\nconst body = context.document.body;\n\nconst results = body.search(\"\uc758 \ubcc0\ub3d9\uc131\uc5d0 \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nresults.items[0].insertText(\"\uc758 \ubcc0\ub3d9\uc131\uc5d0 \", \"Replace\");\nawait context.sync();\n\n// 2) Move the \"_GoBack\" bookmark to the very end of the document, i.e.\n//    right after the last character of the last paragraph and right\n//    before that paragraph's ending mark.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nbody.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) \"...term\uc758 \ubcc0\ub3d9\" + <hidden _GoBack bookmark> + \"\uc131\uc5d0 max(v, 0)...\"\n#    becomes a single, unsplit run \"...term\uc758 \ubcc0\ub3d9\uc131\uc5d0 max(v, 0)...\".\n#    The bookmark that used to sit between the two runs is removed from\n#    here (it gets relocated to the end of the document in step 2).\n# ------------------------------------------------------------------\n\n$goBack = $d.Bookmarks.Item(\"_GoBack\")\n$goBack.Delete()\n\n$mergeRange = $d.Content\n$mergeRange.Find.MatchWildcards = $false\n$mergeRange.Find.Forward = $true\n$mergeRange.Find.Execute(\"\uc758 \ubcc0\ub3d9\uc131\uc5d0 \") | Out-Null\n\n# Re-writing a Range with text identical to what is already there is a\n# no-op, so round-trip through a placeholder to force Word to collapse\n# the (still separately-run) matched text into the single run implied\n# by the assignment.\n$mergeRange.Text = \"IRON_NATIVE_TMP_PLACEHOLDER\"\n\n$mergeRange2 = $d.Content\n$mergeRange2.Find.MatchWildcards = $false\n$mergeRange2.Find.Forward = $true\n$mergeRange2.Find.Execute(\"IRON_NATIVE_TMP_PLACEHOLDER\") | Out-Null\n$mergeRange2.Text = \"\uc758 \ubcc0\ub3d9\uc131\uc5d0 \"\n\n# ------------------------------------------------------------------\n# 2) Move the \"_GoBack\" bookmark to the very end of the document, i.e.\n#    right after the last character of the last paragraph and right\n#    before that paragraph's ending mark.\n# ------------------------------------------------------------------\n\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$endOfText = $lastParagraph.Range.End - 1\n\n$markerRange = $d.Range($endOfText, $endOfText)\n$markerRange.InsertAfter(\"~\")\n# $markerRange now spans exactly the \"~\" placeholder character we just\n# inserted, so it is a safe, non-collapsed anchor for Bookmarks.Add.\n$d.Bookmarks.Add(\"_GoBack\", $markerRange) | Out-Null\n\n$placeholderRange = $d.Range($endOfText, $endOfText + 1)\n$placeholderRange.Delete()\n"}
